{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Target edit (per commit \"#6 Turning point All is set wrapping up today 22/01/2024\"):\n//   - Replace the whole body content with a single paragraph reading\n//     \"Hello, {firstname} {mdname} {sirname}!! You are a programmer!!\"\n//     (dropping the old \"Jay Ganesh!\" line and the {name}/{age}/{hobby}/\n//     {nick name} placeholder paragraphs + trailing empty paragraph).\n//   - Add a \"_GoBack\" bookmark (start/end, no content) at the end of that\n//     paragraph, as Word does automatically after an edit.\n//   - Update the section's page size / header-footer distances from the\n//     \"Letter\" defaults to the \"A4\" defaults that a refreshed template\n//     produced (pgSz 12240x15840 -> 11906x16838, header/footer 720->708).\n//     (Column spacing also moves 720->708 twips in the target, but the\n//     Word JavaScript API has no property for that value, so it is left\n//     at its default here.)\n\nconst body = context.document.body;\n\n// 1) Wipe all existing paragraphs/content; Word always leaves exactly one\n//    empty paragraph behind, which becomes our new single paragraph.\nbody.clear();\nawait context.sync();\n\n// 2) Insert the new sentence into that paragraph.\nbody.insertText(\n  \"Hello, {firstname} {mdname} {sirname}!! You are a programmer!!\",\n  Word.InsertLocation.start\n);\nawait context.sync();\n\n// 3) Word stamps a \"_GoBack\" bookmark at the last edit location; add it at\n//    the end of the (only) paragraph to match.\nconst endRange = body.getRange(Word.RangeLocation.end);\nendRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// 4) Update the page setup (section) to the new A4-based dimensions.\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nconst pageSetup = sections.items[0].pageSetup;\npageSetup.pageWidth = 595.3; // 11906 twips\npageSetup.pageHeight = 841.9; // 16838 twips\npageSetup.headerDistance = 35.4; // 708 twips\npageSetup.footerDistance = 35.4; // 708 twips\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $d ($word.ActiveDocument) is the open document.\n#\n# Target edit (per commit \"#6 Turning point All is set wrapping up today 22/01/2024\"):\n#   - Replace the whole body content with a single paragraph reading\n#     \"Hello, {firstname} {mdname} {sirname}!! You are a programmer!!\"\n#     (dropping the old \"Jay Ganesh!\" line and the {name}/{age}/{hobby}/\n#     {nick name} placeholder paragraphs + trailing empty paragraph).\n#   - Add a \"_GoBack\" bookmark (start/end, no content) right after that\n#     text, as Word does automatically after an edit.\n#   - Update the section's page size / header-footer distances / column\n#     spacing from the \"Letter\" defaults to the \"A4\" defaults that a\n#     refreshed template produced (pgSz 12240x15840 -> 11906x16838,\n#     header/footer 720->708, column spacing 720->708 twips).\n\n$d = $word.ActiveDocument\n\n# 1) Wipe all existing content. Each Delete() on $d.Content removes up to\n#    the next paragraph mark in this host, so re-fetch and repeat until the\n#    whole story is empty; Word always leaves exactly one empty paragraph\n#    behind, which becomes our new single paragraph.\nfor ($i = 0; $i -lt 20; $i++) {\n    $current = $d.Content\n    if ($current.Text -eq \"\") { break }\n    $current.Delete()\n}\n\n# 2) Insert the new sentence, plus a temporary trailing placeholder\n#    character so the bookmark (step 3) can be anchored exactly at the end\n#    of the real text without snapping to the paragraph mark.\n$newText = \"Hello, {firstname} {mdname} {sirname}!! You are a programmer!!\"\n$r = $d.Range()\n$r.InsertAfter($newText + \"X\")\n\n# 3) Word stamps a \"_GoBack\" bookmark at the last edit location; add it\n#    right after the sentence (collapsed, no spanned text) to match.\n$endPos = $newText.Length\n$endRange = $d.Range($endPos, $endPos)\n$d.Bookmarks.Add(\"_GoBack\", $endRange)\n\n# Remove the temporary placeholder character now that the bookmark is set.\n$placeholder = $d.Range($endPos, $endPos + 1)\n$placeholder.Delete()\n\n# 4) Update the page setup (section) to the new A4-based dimensions.\n$ps = $d.PageSetup\n$ps.PageWidth = 595.3        # 11906 twips\n$ps.PageHeight = 841.9       # 16838 twips\n$ps.HeaderDistance = 35.4    # 708 twips\n$ps.FooterDistance = 35.4    # 708 twips\n$ps.TextColumns.Spacing = 35.4   # 708 twips\n"}
